$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Update column A values (rows 2-31) per new parameters
for ($r = 2; $r -le 11; $r++) {
    $ws.Cells.Item($r, 1).Value = 1.3
}
for ($r = 12; $r -le 21; $r++) {
    $ws.Cells.Item($r, 1).Value = 1.5
}
for ($r = 22; $r -le 31; $r++) {
    $ws.Cells.Item($r, 1).Value = 1.7
}

# Update the selection shown in the sheet view
$ws.Activate()
$ws.Range("E2:E4").Select()
